# Scheduled market-data refresh: updates currentAveragePrice(NQ/HQ),
# LevePrice(NQ/HQ) and LeveProfit(NQ/HQ) figures for a batch of Leves
# across the job sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).

$wb = $excel.ActiveWorkbook

function Set-Row($SheetName, $Row, $Values) {
    $ws = $wb.Worksheets.Item($SheetName)
    foreach ($col in $Values.Keys) {
        $ws.Range("$col$Row").Value = $Values[$col]
    }
}

Set-Row "ALC" 18  @{ H=500; I=500; K=500; M=-216 }
Set-Row "ALC" 132 @{ H=4858.6665; I=4858.6665; K=14575.9995; M=-12045.9995 }
Set-Row "ALC" 138 @{ H=2167.4355; I=1110.5; J=3294.8333; K=3331.5; L=9884.499899999999; M=1808.5; N=-20164.4999 }

Set-Row "ARM" 2   @{ H=2257.375; I=1410.6; J=3668.6667; K=1410.6; L=3668.6667; M=-1297.6; N=-3894.6667 }
Set-Row "ARM" 32  @{ H=40033.715; I=43557.2; J=10671.333; K=43557.2; L=10671.333; M=-43270.2; N=-11245.333 }
Set-Row "ARM" 74  @{ H=1735.9166; I=1359.8422; J=3165; K=1359.8422; L=3165; M=-485.8422; N=-4913 }
Set-Row "ARM" 77  @{ H=1735.9166; I=1359.8422; J=3165; K=6799.211; L=15825; M=-2431.211; N=-24561 }
Set-Row "ARM" 114 @{ H=43066.332; J=43066.332; L=43066.332; N=-51744.332 }
Set-Row "ARM" 116 @{ H=2257.375; I=1410.6; J=3668.6667; K=1410.6; L=3668.6667; M=883.4000000000001; N=-8256.6667 }
Set-Row "ARM" 132 @{ H=27734.4; I=1948; J=66414; K=5844; L=199242; M=-3314; N=-204302 }

Set-Row "BSM" 3   @{ H=2257.375; I=1410.6; J=3668.6667; K=1410.6; L=3668.6667; M=-1296.6; N=-3896.6667 }
Set-Row "BSM" 112 @{ H=50000; J=50000; L=50000; N=-52954 }
Set-Row "BSM" 134 @{ H=106080; I=117600; J=2400; K=352800; L=7200; M=-350265; N=-12270 }

Set-Row "CRP" 31  @{ H=10717.122; I=18624.824; J=5115.8335; K=18624.824; L=5115.8335; M=-18329.824; N=-5705.8335 }
Set-Row "CRP" 34  @{ H=10717.122; I=18624.824; J=5115.8335; K=18624.824; L=5115.8335; M=-18422.824; N=-5519.8335 }
Set-Row "CRP" 41  @{ H=13983.333; J=19975; L=19975; N=-20831 }

Set-Row "CUL" 131 @{ H=774.3626; J=774.3626; L=2323.0878; N=-12403.0878 }

Set-Row "GSM" 122 @{ H=1627.091; I=1501.4286; J=1847; K=4504.2858; L=5541; M=-2054.2858; N=-10441 }
Set-Row "GSM" 126 @{ H=3932.2424; I=3079.4119; K=9238.235700000001; M=-6768.235700000001 }
Set-Row "GSM" 132 @{ H=87319.61; I=149046.28; J=48039; K=447138.84; L=144117; M=-444608.84; N=-149177 }

Set-Row "LTW" 22  @{ H=2153.5; J=3160; L=3160; N=-3750 }
Set-Row "LTW" 27  @{ H=2153.5; J=3160; L=3160; N=-3374 }
Set-Row "LTW" 46  @{ H=965.0909; I=723.05884; J=1788; K=723.05884; L=1788; M=-535.05884; N=-2164 }
Set-Row "LTW" 82  @{ H=2335.7144; I=3011.111; J=1120; K=3011.111; L=1120; M=-2650.111; N=-1842 }
Set-Row "LTW" 85  @{ H=2335.7144; I=3011.111; J=1120; K=3011.111; L=1120; M=-1763.111; N=-3616 }
Set-Row "LTW" 100 @{ H=2616.6667; I=1477.7778; J=3755.5557; K=1477.7778; L=3755.5557; M=-936.7778000000001; N=-4837.5557 }
Set-Row "LTW" 110 @{ H=34950; J=34950; L=34950; N=-43130 }

Set-Row "WVR" 136 @{ H=949958.4; I=1537030.5; J=1611.0769; K=4611091.5; L=4833.2307; M=-4608541.5; N=-9933.2307 }
